$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Bump the cached "datetimeFigureOut" field text (11/8/2017 -> 11/13/2018)
#    on the slide master and every custom (slide) layout's Date Placeholder.
# ---------------------------------------------------------------------------
$newDate = "11/13/2018"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Slide 6 ("Quiz 1 Level 4") edits
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)

# Title: the two identically-formatted runs "Quiz 1 Level 4: " and
# "Solve the Following" collapse into a single run with the same text.
$title = $s6.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "ZZZ_TEMP_ZZZ"
$title.TextFrame.TextRange.Text = "Quiz 1 Level 4: Solve the Following"

# Content placeholder: strike through the 3rd paragraph (the "Would it be a
# better or worse idea ..." sentence).
$content = $s6.Shapes.Item(2)
$tr = $content.TextFrame.TextRange
$para = $tr.Paragraphs(3, 1)
$para.Font.Strike = -1
